$d = $word.ActiveDocument

$d.Content.Find.Execute("2025-10-04 Saturday", $true, $false, $false, $false, $false, $true, 1, $false, "2025-10-05 Sunday", 2) | Out-Null
$d.Content.Find.Execute("21-2=19", $true, $false, $false, $false, $false, $true, 1, $false, "34+23=57", 2) | Out-Null
$d.Content.Find.Execute("27+28=55", $true, $false, $false, $false, $false, $true, 1, $false, "55-35=20", 2) | Out-Null
$d.Content.Find.Execute("84-44=40", $true, $false, $false, $false, $false, $true, 1, $false, "16+31=47", 2) | Out-Null
$d.Content.Find.Execute("14+74=88", $true, $false, $false, $false, $false, $true, 1, $false, "52-4=48", 2) | Out-Null
$d.Content.Find.Execute("26-8=18", $true, $false, $false, $false, $false, $true, 1, $false, "99-34=65", 2) | Out-Null
$d.Content.Find.Execute("86-51=35", $true, $false, $false, $false, $false, $true, 1, $false, "33-1=32", 2) | Out-Null
$d.Content.Find.Execute("16+30=46", $true, $false, $false, $false, $false, $true, 1, $false, "7+21=28", 2) | Out-Null
$d.Content.Find.Execute("43-2=41", $true, $false, $false, $false, $false, $true, 1, $false, "79-35=44", 2) | Out-Null
$d.Content.Find.Execute("42+51=93", $true, $false, $false, $false, $false, $true, 1, $false, "53+30=83", 2) | Out-Null
$d.Content.Find.Execute("62+10=72", $true, $false, $false, $false, $false, $true, 1, $false, "31-19=12", 2) | Out-Null
$d.Content.Find.Execute("30-7=23", $true, $false, $false, $false, $false, $true, 1, $false, "1+35=36", 2) | Out-Null
$d.Content.Find.Execute("68+30=98", $true, $false, $false, $false, $false, $true, 1, $false, "73-17=56", 2) | Out-Null
$d.Content.Find.Execute("42+26=68", $true, $false, $false, $false, $false, $true, 1, $false, "97-72=25", 2) | Out-Null
$d.Content.Find.Execute("89-74=15", $true, $false, $false, $false, $false, $true, 1, $false, "46-22=24", 2) | Out-Null
$d.Content.Find.Execute("89-59=30", $true, $false, $false, $false, $false, $true, 1, $false, "23-6=17", 2) | Out-Null
$d.Content.Find.Execute("50-30=20", $true, $false, $false, $false, $false, $true, 1, $false, "2+35=37", 2) | Out-Null
$d.Content.Find.Execute("61+34=95", $true, $false, $false, $false, $false, $true, 1, $false, "45+40=85", 2) | Out-Null
$d.Content.Find.Execute("88-45=43", $true, $false, $false, $false, $false, $true, 1, $false, "25-13=12", 2) | Out-Null
$d.Content.Find.Execute("36-30=6", $true, $false, $false, $false, $false, $true, 1, $false, "67-25=42", 2) | Out-Null
$d.Content.Find.Execute("48+29=77", $true, $false, $false, $false, $false, $true, 1, $false, "43+33=76", 2) | Out-Null
$d.Content.Find.Execute("39-1=38", $true, $false, $false, $false, $false, $true, 1, $false, "2+40=42", 2) | Out-Null
$d.Content.Find.Execute("41+56=97", $true, $false, $false, $false, $false, $true, 1, $false, "67+9=76", 2) | Out-Null
$d.Content.Find.Execute("81-59=22", $true, $false, $false, $false, $false, $true, 1, $false, "68-34=34", 2) | Out-Null
$d.Content.Find.Execute("10+73=83", $true, $false, $false, $false, $false, $true, 1, $false, "46-16=30", 2) | Out-Null
$d.Content.Find.Execute("45-5=40", $true, $false, $false, $false, $false, $true, 1, $false, "55-6=49", 2) | Out-Null
$d.Content.Find.Execute("75-58=17", $true, $false, $false, $false, $false, $true, 1, $false, "43+25=68", 2) | Out-Null
$d.Content.Find.Execute("58+40=98", $true, $false, $false, $false, $false, $true, 1, $false, "56-9=47", 2) | Out-Null
$d.Content.Find.Execute("49+3=52", $true, $false, $false, $false, $false, $true, 1, $false, "69+27=96", 2) | Out-Null
$d.Content.Find.Execute("52+10=62", $true, $false, $false, $false, $false, $true, 1, $false, "99-7=92", 2) | Out-Null
$d.Content.Find.Execute("6+35=41", $true, $false, $false, $false, $false, $true, 1, $false, "51+23=74", 2) | Out-Null
$d.Content.Find.Execute("51+4=55", $true, $false, $false, $false, $false, $true, 1, $false, "45-4=41", 2) | Out-Null
$d.Content.Find.Execute("62+5=67", $true, $false, $false, $false, $false, $true, 1, $false, "1+46=47", 2) | Out-Null
$d.Content.Find.Execute("85-64=21", $true, $false, $false, $false, $false, $true, 1, $false, "30+33=63", 2) | Out-Null
$d.Content.Find.Execute("44+40=84", $true, $false, $false, $false, $false, $true, 1, $false, "99-73=26", 2) | Out-Null
$d.Content.Find.Execute("94-28=66", $true, $false, $false, $false, $false, $true, 1, $false, "3+60=63", 2) | Out-Null
$d.Content.Find.Execute("37-11=26", $true, $false, $false, $false, $false, $true, 1, $false, "89-19=70", 2) | Out-Null
$d.Content.Find.Execute("65+29=94", $true, $false, $false, $false, $false, $true, 1, $false, "5+43=48", 2) | Out-Null
$d.Content.Find.Execute("46+47=93", $true, $false, $false, $false, $false, $true, 1, $false, "97-47=50", 2) | Out-Null
$d.Content.Find.Execute("55-14=41", $true, $false, $false, $false, $false, $true, 1, $false, "26-23=3", 2) | Out-Null
$d.Content.Find.Execute("54-33=21", $true, $false, $false, $false, $false, $true, 1, $false, "86-22=64", 2) | Out-Null
$d.Content.Find.Execute("52-0=52", $true, $false, $false, $false, $false, $true, 1, $false, "4+54=58", 2) | Out-Null
$d.Content.Find.Execute("45-44=1", $true, $false, $false, $false, $false, $true, 1, $false, "51-5=46", 2) | Out-Null
$d.Content.Find.Execute("45+10=55", $true, $false, $false, $false, $false, $true, 1, $false, "81-81=0", 2) | Out-Null
$d.Content.Find.Execute("28+15=43", $true, $false, $false, $false, $false, $true, 1, $false, "87-21=66", 2) | Out-Null
$d.Content.Find.Execute("40-12=28", $true, $false, $false, $false, $false, $true, 1, $false, "92-70=22", 2) | Out-Null
$d.Content.Find.Execute("43+7=50", $true, $false, $false, $false, $false, $true, 1, $false, "88-27=61", 2) | Out-Null
$d.Content.Find.Execute("79-6=73", $true, $false, $false, $false, $false, $true, 1, $false, "12+38=50", 2) | Out-Null
$d.Content.Find.Execute("96-52=44", $true, $false, $false, $false, $false, $true, 1, $false, "93-55=38", 2) | Out-Null
$d.Content.Find.Execute("11+73=84", $true, $false, $false, $false, $false, $true, 1, $false, "86-20=66", 2) | Out-Null
$d.Content.Find.Execute("83-51=32", $true, $false, $false, $false, $false, $true, 1, $false, "92-15=77", 2) | Out-Null
$d.Content.Find.Execute("94-86=8", $true, $false, $false, $false, $false, $true, 1, $false, "41+35=76", 2) | Out-Null
$d.Content.Find.Execute("30-9=21", $true, $false, $false, $false, $false, $true, 1, $false, "97-42=55", 2) | Out-Null
$d.Content.Find.Execute("80-17=63", $true, $false, $false, $false, $false, $true, 1, $false, "53-16=37", 2) | Out-Null
$d.Content.Find.Execute("90-9=81", $true, $false, $false, $false, $false, $true, 1, $false, "31-8=23", 2) | Out-Null
$d.Content.Find.Execute("89+7=96", $true, $false, $false, $false, $false, $true, 1, $false, "10+82=92", 2) | Out-Null
$d.Content.Find.Execute("59+33=92", $true, $false, $false, $false, $false, $true, 1, $false, "47-23=24", 2) | Out-Null
$d.Content.Find.Execute("58+27=85", $true, $false, $false, $false, $false, $true, 1, $false, "19-15=4", 2) | Out-Null
$d.Content.Find.Execute("97-43=54", $true, $false, $false, $false, $false, $true, 1, $false, "84-56=28", 2) | Out-Null
$d.Content.Find.Execute("52+45=97", $true, $false, $false, $false, $false, $true, 1, $false, "53-9=44", 2) | Out-Null
$d.Content.Find.Execute("70-43=27", $true, $false, $false, $false, $false, $true, 1, $false, "15+42=57", 2) | Out-Null
$d.Content.Find.Execute("33-12=21", $true, $false, $false, $false, $false, $true, 1, $false, "72-13=59", 2) | Out-Null
$d.Content.Find.Execute("78-67=11", $true, $false, $false, $false, $false, $true, 1, $false, "98-1=97", 2) | Out-Null
$d.Content.Find.Execute("87+0=87", $true, $false, $false, $false, $false, $true, 1, $false, "78+17=95", 2) | Out-Null
$d.Content.Find.Execute("58-7=51", $true, $false, $false, $false, $false, $true, 1, $false, "92-79=13", 2) | Out-Null
$d.Content.Find.Execute("91-65=26", $true, $false, $false, $false, $false, $true, 1, $false, "18+70=88", 2) | Out-Null
$d.Content.Find.Execute("84-75=9", $true, $false, $false, $false, $false, $true, 1, $false, "81-25=56", 2) | Out-Null
$d.Content.Find.Execute("93-36=57", $true, $false, $false, $false, $false, $true, 1, $false, "73-20=53", 2) | Out-Null
$d.Content.Find.Execute("89-73=16", $true, $false, $false, $false, $false, $true, 1, $false, "21+77=98", 2) | Out-Null
$d.Content.Find.Execute("0+19=19", $true, $false, $false, $false, $false, $true, 1, $false, "11+36=47", 2) | Out-Null
$d.Content.Find.Execute("19+15=34", $true, $false, $false, $false, $false, $true, 1, $false, "81-47=34", 2) | Out-Null
$d.Content.Find.Execute("48+42=90", $true, $false, $false, $false, $false, $true, 1, $false, "78+18=96", 2) | Out-Null
$d.Content.Find.Execute("43-40=3", $true, $false, $false, $false, $false, $true, 1, $false, "46-9=37", 2) | Out-Null
$d.Content.Find.Execute("47+5=52", $true, $false, $false, $false, $false, $true, 1, $false, "30+42=72", 2) | Out-Null
$d.Content.Find.Execute("9+13=22", $true, $false, $false, $false, $false, $true, 1, $false, "3+72=75", 2) | Out-Null
$d.Content.Find.Execute("79-28=51", $true, $false, $false, $false, $false, $true, 1, $false, "43+46=89", 2) | Out-Null
$d.Content.Find.Execute("49-2=47", $true, $false, $false, $false, $false, $true, 1, $false, "7+56=63", 2) | Out-Null
$d.Content.Find.Execute("46+7=53", $true, $false, $false, $false, $false, $true, 1, $false, "35-7=28", 2) | Out-Null
$d.Content.Find.Execute("42+54=96", $true, $false, $false, $false, $false, $true, 1, $false, "24+67=91", 2) | Out-Null
$d.Content.Find.Execute("8+31=39", $true, $false, $false, $false, $false, $true, 1, $false, "92-89=3", 2) | Out-Null
$d.Content.Find.Execute("53-17=36", $true, $false, $false, $false, $false, $true, 1, $false, "40+50=90", 2) | Out-Null
$d.Content.Find.Execute("93-52=41", $true, $false, $false, $false, $false, $true, 1, $false, "4+77=81", 2) | Out-Null
$d.Content.Find.Execute("81-54=27", $true, $false, $false, $false, $false, $true, 1, $false, "90-89=1", 2) | Out-Null
$d.Content.Find.Execute("25+27=52", $true, $false, $false, $false, $false, $true, 1, $false, "88-73=15", 2) | Out-Null
$d.Content.Find.Execute("78-0=78", $true, $false, $false, $false, $false, $true, 1, $false, "54+18=72", 2) | Out-Null
$d.Content.Find.Execute("90-32=58", $true, $false, $false, $false, $false, $true, 1, $false, "23-19=4", 2) | Out-Null
$d.Content.Find.Execute("56-53=3", $true, $false, $false, $false, $false, $true, 1, $false, "91-48=43", 2) | Out-Null
$d.Content.Find.Execute("11+56=67", $true, $false, $false, $false, $false, $true, 1, $false, "49-24=25", 2) | Out-Null
$d.Content.Find.Execute("84-27=57", $true, $false, $false, $false, $false, $true, 1, $false, "74+12=86", 2) | Out-Null
$d.Content.Find.Execute("53-13=40", $true, $false, $false, $false, $false, $true, 1, $false, "29+15=44", 2) | Out-Null
$d.Content.Find.Execute("96-65=31", $true, $false, $false, $false, $false, $true, 1, $false, "33+37=70", 2) | Out-Null
$d.Content.Find.Execute("16+61=77", $true, $false, $false, $false, $false, $true, 1, $false, "42-3=39", 2) | Out-Null
$d.Content.Find.Execute("97-10=87", $true, $false, $false, $false, $false, $true, 1, $false, "17+56=73", 2) | Out-Null
$d.Content.Find.Execute("49+7=56", $true, $false, $false, $false, $false, $true, 1, $false, "64+27=91", 2) | Out-Null
$d.Content.Find.Execute("3+81=84", $true, $false, $false, $false, $false, $true, 1, $false, "29+17=46", 2) | Out-Null
$d.Content.Find.Execute("54-3=51", $true, $false, $false, $false, $false, $true, 1, $false, "59-49=10", 2) | Out-Null
$d.Content.Find.Execute("95-85=10", $true, $false, $false, $false, $false, $true, 1, $false, "17+50=67", 2) | Out-Null
$d.Content.Find.Execute("45+42=87", $true, $false, $false, $false, $false, $true, 1, $false, "57-48=9", 2) | Out-Null
$d.Content.Find.Execute("38-19=19", $true, $false, $false, $false, $false, $true, 1, $false, "82-56=26", 2) | Out-Null
$d.Content.Find.Execute("32+19=51", $true, $false, $false, $false, $false, $true, 1, $false, "86-50=36", 2) | Out-Null
$d.Content.Find.Execute("20+1=21", $true, $false, $false, $false, $false, $true, 1, $false, "63+21=84", 2) | Out-Null
